# Deep sea double count fix
# Updates computed weight-percentage values for row 7 (Gadus morhua)
# and row 14 (Global) after correcting a double-counting bug in the
# upstream aggregation for deep sea species.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - Gadus morhua
$ws.Range("C7").Value = 0.684102635418281
$ws.Range("D7").Value = 0.4612462745817191
$ws.Range("F7").Value = 59.72875422025597
$ws.Range("G7").Value = 40.27124577974401

# Row 14 - Global
$ws.Range("C14").Value = 11.73967789285522
$ws.Range("D14").Value = 3.130923048290898
$ws.Range("F14").Value = 52.70387400918387
$ws.Range("G14").Value = 14.05590301331859
$ws.Range("H14").Value = 33.24022297749755
